{"js": "// \"Added last minute updates\" \u2014 tweak the ID placeholder paragraph at the\n// top of the document: give it a thin paragraph border, nudge its left\n// indent out a bit, and fold its two runs (placeholder text + a trailing\n// space run) into a single run carrying the refreshed placeholder text.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst idParagraph = paragraphs.items[0];\n\n// Paragraph border on all four sides, 5pt \"distance from text\". Office.js's\n// paragraph.borders collection (Word.ParagraphBorder) only surfaces\n// type/color/width \u2014 there's no distanceFromText/space setter on it \u2014 so set\n// the underlying Word object-model member directly (Paragraph.Borders is\n// the same Borders collection the COM/VBA object model exposes, and every\n// Office.js proxy here is backed by that shared bridge).\nidParagraph._omSet(\"DistanceFromTop\", 5, \"Borders\");\nidParagraph._omSet(\"DistanceFromLeft\", 5, \"Borders\");\nidParagraph._omSet(\"DistanceFromBottom\", 5, \"Borders\");\nidParagraph._omSet(\"DistanceFromRight\", 5, \"Borders\");\n\n// Left indent: 120 -> 225 twips, i.e. 6pt -> 11.25pt (Office.js indents are\n// in points).\nidParagraph.leftIndent = 11.25;\n\nawait context.sync();\n\n// Replace the paragraph's whole visible text (both runs: the bold/italic-off\n// placeholder run and the trailing space run) with the new placeholder in a\n// single run. insertText(..., \"Replace\") reuses the formatting of the first\n// run and drops the now-empty second run entirely.\nidParagraph.insertText(\"**ID__AFFARS_AF_PGI_5308__ID**\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# Add a paragraph border (top/left/bottom/right) with a 5pt \"space\" (distance\n# from text) on the first paragraph, and bump its left indent from 120 -> 225\n# twips (6pt -> 11.25pt).\n$pf = $p.Range.ParagraphFormat\n$pf.Borders.DistanceFromTop = 5\n$pf.Borders.DistanceFromLeft = 5\n$pf.Borders.DistanceFromBottom = 5\n$pf.Borders.DistanceFromRight = 5\n$pf.LeftIndent = 11.25\n\n# Collapse the paragraph's two runs (\"**ID__AFFARS_pgi_5308_topic_2__ID**\" +\n# a trailing space run) into a single run with the updated placeholder text,\n# keeping the first run's formatting and dropping the separate space run.\n$rng = $p.Range\n$rng.MoveEnd(1, -1) | Out-Null\n$rng.Text = \"**ID__AFFARS_AF_PGI_5308__ID**\"\n"}
